# Bulk match upload test fixture update
# - Adds "league" to the Match Type header hint
# - Backfills missing Player DOB values (columns E/H) for rows 3 and 4
# - Shifts the previously mis-aligned Is-Doubles/Score/Location/Notes values
#   in rows 3 and 4 one column to the right so they line up with the header
# - Changes row 5's match type to "league" and updates its notes text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value = "Match Type (casual/tournament/league)"

# --- Row 3 ---
# Add missing Player 1 / Player 2 date-of-birth values
$ws.Range("E3").Value = 31116
$ws.Range("H3").Value = 33076

# Re-align the trailing columns (Is Doubles .. Notes) one column to the right.
# Write the new (shifted) positions first using the values that used to live
# one column to the left, then clear whatever is left over.
$ws.Range("O3").Value = $false
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 11
$ws.Range("S3").Value = 8
$ws.Range("T3").Value = 11
$ws.Range("U3").Value = 6
$ws.Range("X3").Value = "Tournament Court A"
$ws.Range("Y3").Value = "Singles tournament match - straight sets"

# Clear the now-stale cells that are no longer part of row 3's data
$ws.Range("N3").ClearContents()
$ws.Range("W3").ClearContents()

# --- Row 4 ---
# Add missing Player 1 / Player 2 date-of-birth values
$ws.Range("E4").Value = 32452
$ws.Range("H4").Value = 33865

# Re-align the trailing columns (Is Doubles .. Notes) one column to the right.
$ws.Range("O4").Value = $false
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 9
$ws.Range("S4").Value = 11
$ws.Range("T4").Value = 10
$ws.Range("U4").Value = 12
$ws.Range("V4").Value = 11
$ws.Range("W4").Value = 13
$ws.Range("X4").Value = "Court 2"
$ws.Range("Y4").Value = "Close singles match with three games"

# Clear the now-stale cell that is no longer part of row 4's data
$ws.Range("N4").ClearContents()

# --- Row 5 ---
$ws.Range("B5").Value = "league"
$ws.Range("Y5").Value = "League doubles match - Team 2 wins in straight sets"
